# Generate Report for Handoff
# Updates the localization-status report for the e4599166 file:
# it has moved from "Handed back: in sync with en-US" to "Ready for handoff",
# with a refreshed handoff timestamp and a stale-handback-version error message.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"

$newOverviewDate = "2016-08-19 16:52:05"

$newZhHandoffDate = "2016-08-19 16:51:57"
$newDeHandoffDate = "2016-08-19 16:52:05"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb87d40792e215ccbccf7edaf4aef579327f417a/e2e/e4599166-965f-4bb1-affa-3efa79602bf8.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b710b485cf0fa4da16988fc5d751a19952a3531/e2e/e4599166-965f-4bb1-affa-3efa79602bf8.md."

# --- Overview sheet: row 3 is the e4599166 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = $newOverviewDate

# --- zh-cn sheet: row 3 is the e4599166 file ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusNew
$wsZh.Range("H3").Value = $newZhHandoffDate
$wsZh.Range("P3").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row 3 is the e4599166 file ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusNew
$wsDe.Range("H3").Value = $newDeHandoffDate
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.17
